$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the B1 header: " atf_analysis_id" (typo, leading space) -> "matf_analysis_id"
$ws.Range("B1").Value = "matf_analysis_id"

# Update the data validation on column C to describe the Data_file table ID
# (was "Data_file_name" referencing the SRA filename field)
$val = $ws.Range("C1").Validation
$val.InputTitle = "ID from Data_file table"
$val.InputMessage = "This ID has to already exist in the Data_file table in the Data_file_id"

# Restore the selection to column C (as last selected/active range)
$ws.Range("C1:C1048576").Select()
